$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sector / average-correlation data for rows 3 through 53.
# Row 1 is the header and row 2 (Industrial Conglomerates) is unchanged.
$data = @(
    @(3, 'Road & Rail(22)', 0.5989653154570562),
    @(4, 'Air Freight & Logistics(11)', 0.5775675828705409),
    @(5, 'Energy Equipment & Services(32)', 0.5751299603686909),
    @(6, 'Marine(15)', 0.5690460143460049),
    @(7, 'Construction & Engineering(20)', 0.5675362280049581),
    @(8, 'Construction Materials(8)', 0.5663544204918297),
    @(9, 'Trading Companies & Distributors(25)', 0.5441324810134711),
    @(10, 'Containers & Packaging(12)', 0.5275944728785374),
    @(11, 'Machinery(85)', 0.5214815206989146),
    @(12, 'Building Products(23)', 0.5123627667862571),
    @(13, 'Metals & Mining(89)', 0.5109642925283362),
    @(14, 'Auto Components(21)', 0.5046759631775736),
    @(15, 'Multi-Utilities(18)', 0.4909455662030818),
    @(16, 'Chemicals(51)', 0.4782659817295227),
    @(17, 'Life Sciences Tools & Services(19)', 0.4662222783206468),
    @(18, 'Wireless Telecommunication Services(14)', 0.4456355275386176),
    @(19, 'Airlines(14)', 0.4420316194123448),
    @(20, 'Insurance(75)', 0.4331316156024389),
    @(21, 'Gas Utilities(12)', 0.4197610392852889),
    @(22, 'Capital Markets(75)', 0.4179694322295037),
    @(23, 'Semiconductors & Semiconductor Equipment(68)', 0.4163578740925429),
    @(24, 'IT Services(52)', 0.4142540919728485),
    @(25, 'Leisure Products(11)', 0.4131002866899252),
    @(26, 'Electrical Equipment(28)', 0.4110709898800927),
    @(27, 'Oil, Gas & Consumable Fuels(122)', 0.4037955009463311),
    @(28, 'Household Durables(39)', 0.4036972876281053),
    @(29, 'Professional Services(35)', 0.3793475327116087),
    @(30, 'Water Utilities(12)', 0.3783444804567622),
    @(31, 'Health Care Providers & Services(46)', 0.3740882294784622),
    @(32, 'Electric Utilities(28)', 0.3626800602595638),
    @(33, 'Communications Equipment(45)', 0.3451938496260053),
    @(34, 'Banks(246)', 0.3416126992819677),
    @(35, 'Consumer Finance(15)', 0.3360127359773946),
    @(36, 'Specialty Retail(58)', 0.3341948423933782),
    @(37, 'Food & Staples Retailing(15)', 0.3328140188231028),
    @(38, 'Aerospace & Defense(37)', 0.3327326980024479),
    @(39, 'Software(66)', 0.3316468995720376),
    @(40, 'Hotels, Restaurants & Leisure(50)', 0.3270062579011052),
    @(41, 'Commercial Services & Supplies(52)', 0.3159174398305094),
    @(42, 'Textiles, Apparel & Luxury Goods(29)', 0.3084588888265488),
    @(43, 'Beverages(21)', 0.3010316192346816),
    @(44, 'Diversified Consumer Services(17)', 0.2979870840024118),
    @(45, 'Real Estate Management & Development(22)', 0.2665444198014706),
    @(46, 'Entertainment(22)', 0.2468887503039209),
    @(47, 'Media(42)', 0.2466866844775422),
    @(48, 'Diversified Telecommunication Services(20)', 0.2287517233842288),
    @(49, 'Health Care Equipment & Supplies(83)', 0.2180506714605436),
    @(50, 'Food Products(44)', 0.1978768024882402),
    @(51, 'Thrifts & Mortgage Finance(47)', 0.1785330825979096),
    @(52, 'Biotechnology(126)', 0.1760245794419086),
    @(53, 'Pharmaceuticals(48)', 0.1356876458472174)
)

foreach ($item in $data) {
    $r = $item[0]
    $sector = $item[1]
    $value = $item[2]
    $ws.Cells.Item($r, 1).Value = $sector
    $ws.Cells.Item($r, 2).Value = $value
}

# Rows 54-56 no longer exist after the edit (three sectors were folded
# into "no sector" and dropped), so clear their previous contents.
$ws.Range("A54:B56").ClearContents()
